$d = $word.ActiveDocument

# --- Image 1 (InlineShape #5, rId8 / image5.png): editId 0EC75802 -> 3C89FE56 ---
# wp:extent/a:ext cx 5943600 -> 5934075 EMU (467.25pt); cy unchanged (2893695 EMU / 227.85pt)
$shp1 = $d.InlineShapes.Item(5)
$shp1.LockAspectRatio = $false
$shp1.Width = 467.25
$shp1.Height = 227.85
$shp1.LockAspectRatio = $true

# a:srcRect t="5003" -> l="160" t="5003" r="1"
$pf1 = $shp1.PictureFormat
$pf1.CropLeft = 1.585421298389567
$pf1.CropRight = 0.009908883114934795

# --- Image 2 (InlineShape #6, rId9 / image6.png): editId 24C22EB9 -> 74745BD4 ---
# wp:extent/a:ext cx 5943600 -> 5934075 EMU (467.25pt); cy unchanged (2928620 EMU / 230.6pt)
$shp2 = $d.InlineShapes.Item(6)
$shp2.LockAspectRatio = $false
$shp2.Width = 467.25
$shp2.Height = 230.6
$shp2.LockAspectRatio = $true

# a:srcRect t="4947" -> l="160" t="4947" r="1"
$pf2 = $shp2.PictureFormat
$pf2.CropLeft = 1.569819120585582
$pf2.CropRight = 0.009811369503659888

Write-Output "Resized and re-cropped both screenshots (5943600 -> 5934075 EMU width)"
